$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---- Title ----
Replace-Text "Domesticating Space: A Cosmic Leap in Human Frontiers" "A Deeper Insight into the World of Chemistry: Unveiling the Invisible"

# ---- Author ----
Replace-Text "Harriet McCallister" "Dr. Sarah Kingsley"

# ---- Email (leave the trailing '.edu' runs untouched) ----
Replace-Text "harriette.mccallister@spaceacademy" "sarahkingsley@highschool"

# ---- Body paragraph (space essay -> chemistry essay) ----
Replace-Text "Humankind has always looked up at the celestial tapestry above, drawn by its enigmatic beauty and daunting vastness" "The world around us is a complex symphony of elements, molecules, and reactions that orchestrate the very fabric of existence"

Replace-Text " Early civilizations used stars for divining fate, marking seasons, and navigating intricate earthly journeys" " Chemistry, the science that delves into these minute interactions, unveils the secrets hidden within the tangible and intangible, enriching our understanding of the universe we inhabit"

Replace-Text " Today, the cosmos still beckons us, and our ingenuity has provided novel means to explore it" " Embarking on this exploration, we will uncover the fundamental principles that govern the behavior of matter, unveiling the enigmatic dance of elements as they transform into new substances"

# Two sentences (plus the period joining them) collapse into a single new sentence.
Replace-Text " The domestication of space is not a distant aspiration; it is an ongoing endeavor that will reshape our understanding of ourselves and our place in the universe" " As we unravel the mysteries of chemistry, we unlock the potential to harness its powers for the betterment of humanity, pushing the boundaries of scientific knowledge and technological innovation"
Replace-Text ". This essay delves into the scientific, technological, and ethical implications of our ongoing efforts to transform the vast emptiness beyond Earth into a habitable space." "."

Replace-Text "Our conquest of space has commenced with satellites, spacecraft, and space stations--artificial satellites that encircle the Earth, providing indispensable services like communication, global positioning, and weather forecasting" "In the vast laboratory of nature, chemistry orchestrates a mesmerizing symphony of life-sustaining reactions"

Replace-Text " These achievements, however, merely mark the tip of the extraterrestrial iceberg" " From the intricate choreography of photosynthesis to the alchemy of digestion, the very essence of existence is underpinned by chemical processes"

Replace-Text " We aspire for more--permanent habitation in outer space, space tourism, mining lunar and asteroid resources, and ultimately, venturing beyond our solar system" " Chemistry empowers us to unravel these mysteries, deciphering the language of molecules and unlocking the secrets of life itself"

Replace-Text " Yet, to manifest these futuristic dreams, we must overcome formidable challenges: creating artificial gravity, shielding from the harsh radiation of space, developing closed-loop life-support systems, and addressing the psychological tolls of isolation on astronauts" " By comprehending the chemistry of living organisms, we gain invaluable insights into the remarkable resilience and adaptability of life on Earth, inspiring us to protect and preserve the delicate equilibrium of our planet"

Replace-Text "As our presence in space expands, so do the ethical considerations" "Chemistry transcends the confines of the laboratory, shaping our everyday experiences in profound ways"

Replace-Text " Colonizing other celestial bodies raises voprosy of ownership and planetary rights" " From the tantalizing aromas of our favorite foods to the intricate mechanisms that power our technologies, chemistry is an omnipresent force that impacts every facet of our lives"

# Last sentence of that paragraph turns into two sentences (new period added).
Replace-Text " How can we prevent the commoditization and exploitation of extraterrestrial resources? How do we protect and preserve pristine environments like Mars from earthly contamination? What are our responsibilities towards undiscovered life-forms that might exist on other worlds? These questions challenge us to consider our obligations as a species venturing into the great unknown" " It empowers us to create innovative materials, develop life-saving medicines, and safeguard the environment, ensuring a sustainable future for generations to come. By illuminating the fundamental principles of chemistry, we unveil the hidden forces that shape our world, empowering us to become informed and responsible citizens, capable of navigating the complexities of a chemistry-driven society"

# ---- Summary body paragraph ----
Replace-Text "The domestication of space is a daunting undertaking that requires us to push the boundaries of science, technology, and ethics" "This essay delved into the captivating world of chemistry, unveiling the fundamental principles that govern the behavior of matter and highlighting its profound impact on life and society"

Replace-Text " While the challenges are substantial, the potential rewards are immense" " Chemistry empowers us to decipher the mysteries of the universe, harness its powers for the betterment of humanity, and understand the intricate mechanisms that shape our existence"

# Three sentences collapse into one.
Replace-Text " As we embark on this cosmic journey, we must strive for sustainability, equity, and international cooperation" " As we continue to explore the vast tapestry of chemistry, we unlock the potential for scientific breakthroughs, technological advancements, and a deeper appreciation for the wonders of the world around us"
Replace-Text ". The future of humankind may very well lie not only on Earth but among the celestial bodies that have long held our fascination. Our success in domesticating space will be a testament to our ingenuity, adaptability, and our enduring spirit of exploration." "."

# ---- Add a trailing empty paragraph at the very end of the document ----
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
